$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for optimality/feasibility cut counts, matching the
# existing header style (bold/centered/bordered) used by A1:E1.
$ws.Range("F1").Value = "optimality_cuts"
$ws.Range("G1").Value = "feasibility_cuts"
$ws.Range("A1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)

# Update the existing data row (row 2) with refreshed convergence values
# and the new cut-count columns.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = 384399.8731553834
$ws.Range("C2").Value = 424569.0608414363
$ws.Range("D2").Value = 40169.18768605293
$ws.Range("E2").Value = 0.09461166955134041
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5

# Add the two further iterations reported by the benders convergence log.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 384399.8731553834
$ws.Range("C3").Value = 424569.0608414363
$ws.Range("D3").Value = 40169.18768605293
$ws.Range("E3").Value = 0.09461166955134041
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 5

$ws.Range("A4").Value = 8
$ws.Range("B4").Value = 384399.8731553834
$ws.Range("C4").Value = 424569.0608414363
$ws.Range("D4").Value = 40169.18768605293
$ws.Range("E4").Value = 0.09461166955134041
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 5
